# Applies the edit described in the commit:
# "Added method to display brand values as core entities"
#
# 1. Populate column C ("area/point") for a number of rows that previously
#    had no value in that column - most get "area", one gets "area/point".
# 2. Change the separator used in one cell's text from comma to pipe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that should get the value "area" in column C
$areaRows = @(
    308,
    327,328,329,330,331,332,333,334,335,336,337,338,339,340,341,342,
    344,345,
    352,
    386,387,388,
    393,394,395,396,
    398,399,400,
    404,
    407,408,409,
    412,413,414,415,416,417,
    421,
    423,
    430,
    432,
    436,
    452,453,454,
    460,461,462,
    474,
    489,
    491,
    493,494,495,496,497,498,499,500,501,502,503,504,505,506,507,508,509,510,511,512,513,514,515,516,517,518,519,520,521,
    524,525,526,527,
    529,
    536,
    538,539,
    542,543,544,545,546,547,548,549,550,551,552,553,
    567,
    578,579,
    581,582,583,
    586,
    589,590,591
)

foreach ($r in $areaRows) {
    $ws.Range("C$r").Value = "area"
}

# Row 353 gets "area/point" instead of plain "area"
$ws.Range("C353").Value = "area/point"

# Change the pipe-delimited synonyms cell (comma separated -> pipe separated)
$ws.Range("B426").Value = "name|brand|brand name"

# Leave the user's cursor where the edits were made, matching the saved
# view state (active selection on the last touched cell).
$ws.Range("C593").Select()
